$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8
$ws.Range("F8").Value = 2.64
$ws.Range("G8").Value = 2.86
$ws.Range("H8").Value = 2.56
$ws.Range("I8").Value = 2.72
$ws.Range("J8").Value = 3.65
$ws.Range("P8").Value = 2.22
$ws.Range("Q8").Value = 1.69

# Row 9
$ws.Range("F9").Value = 2.06

# Row 10
$ws.Range("N10").Value = 3.95
$ws.Range("O10").Value = 1.31
$ws.Range("Q10").Value = 1.94
$ws.Range("Z10").Value = 10.5

# Row 14
$ws.Range("G14").Value = 2.24

# Row 16
$ws.Range("F16").Value = 2.24
$ws.Range("G16").Value = 2.42
$ws.Range("H16").Value = 3.4
$ws.Range("I16").Value = 3.8
$ws.Range("J16").Value = 3.25
$ws.Range("K16").Value = 3.6
